# Table_S6_RTqPCR-dCT_ms.xlsx — rename the "Mu promoter"/"gene promoter"
# category labels in column C to "Mu TSS"/"gene TSS" (transcription start
# site, not promoter), shrink column C to fit the now-shorter text, nudge
# the print scale, and leave the last worked cell selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Transcript") holds either "Mu promoter" or "gene promoter" for
# every data row (2-35). Rewrite each cell in place so Excel's shared-string
# table drops the two retired strings and appends the renamed ones.
for ($r = 2; $r -le 35; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value()
    if ($val -eq "Mu promoter") {
        $cell.Value = "Mu TSS"
    } elseif ($val -eq "gene promoter") {
        $cell.Value = "gene TSS"
    }
}

# The new labels are shorter, so the best-fit column width shrinks too.
$ws.Columns.Item(3).ColumnWidth = 10.17

# Bump the print scale slightly.
$ws.PageSetup.Zoom = 95

# Leave the view scrolled down with E21 as the active selection.
$ws.Range("E21").Select() | Out-Null
